$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 10.47960307264713
$ws.Range("D2").Value = 3.838247384597163
$ws.Range("E2").Value = 12.7676671299662
$ws.Range("F2").Value = 24.38713471282616
$ws.Range("G2").Value = 30.25043463407922
$ws.Range("H2").Value = 13.72114037120489
$ws.Range("I2").Value = 22.53310231119873
$ws.Range("L2").Value = 9.329841130967058
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("O2").Value = 21.28330997490936

# Row 3
$ws.Range("C3").Value = 10.45037202107951
$ws.Range("D3").Value = 3.835930178954344
$ws.Range("E3").Value = 12.72644113586969
$ws.Range("F3").Value = 24.05226021884549
$ws.Range("G3").Value = 29.53861713261295
$ws.Range("H3").Value = 13.68306165960474
$ws.Range("I3").Value = 22.35737378144339
$ws.Range("L3").Value = 9.328454093328636
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("O3").Value = 21.09573507956754

# Row 4
$ws.Range("C4").Value = 10.43465040749201
$ws.Range("D4").Value = 3.834424825415732
$ws.Range("E4").Value = 12.70398083294624
$ws.Range("F4").Value = 23.85142729026402
$ws.Range("G4").Value = 29.10242667723967
$ws.Range("H4").Value = 13.66258951676806
$ws.Range("I4").Value = 22.25470344097576
$ws.Range("L4").Value = 9.329295013486748
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("O4").Value = 20.9853999308325

# Row 5
$ws.Range("C5").Value = 10.42880791258619
$ws.Range("D5").Value = 3.833790664373051
$ws.Range("E5").Value = 12.69555162918332
$ws.Range("F5").Value = 23.77089341594922
$ws.Range("G5").Value = 28.92518720548097
$ws.Range("H5").Value = 13.65498394736395
$ws.Range("I5").Value = 22.2142189993694
$ws.Range("L5").Value = 9.330063813335171
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 20.94169908310518

# Row 6
$ws.Range("C6").Value = 10.42787196222153
$ws.Range("D6").Value = 3.833684112561818
$ws.Range("E6").Value = 12.69419583309908
$ws.Range("F6").Value = 23.75760278919321
$ws.Range("G6").Value = 28.89579634248972
$ws.Range("H6").Value = 13.6537657119091
$ws.Range("I6").Value = 22.20757951812507
$ws.Range("L6").Value = 9.330217219497447
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 20.93452006577364

# Row 7
$ws.Range("C7").Value = 10.43456932385955
$ws.Range("D7").Value = 3.834416356728021
$ws.Range("E7").Value = 12.70386421649469
$ws.Range("F7").Value = 23.85033575306436
$ws.Range("G7").Value = 29.10003389911527
$ws.Range("H7").Value = 13.66248395444727
$ws.Range("I7").Value = 22.25415191823602
$ws.Range("L7").Value = 9.329303655928605
$ws.Range("N7").Value = 18.02277304767603
$ws.Range("O7").Value = 20.98480539963678

# Row 8
$ws.Range("C8").Value = 10.46906490508427
$ws.Range("D8").Value = 3.83746557992337
$ws.Range("E8").Value = 12.75286433203494
$ws.Range("F8").Value = 24.27074128271646
$ws.Range("G8").Value = 30.00499325919393
$ws.Range("H8").Value = 13.70741082898358
$ws.Range("I8").Value = 22.47145089682981
$ws.Range("L8").Value = 9.329011947853408
$ws.Range("N8").Value = 18.79364780656867
$ws.Range("O8").Value = 21.21765833673262

# Row 9
$ws.Range("C9").Value = 10.55415950457142
$ws.Range("D9").Value = 3.84278785579802
$ws.Range("E9").Value = 12.87127721233435
$ws.Range("F9").Value = 25.12838031552727
$ws.Range("G9").Value = 31.77406878043176
$ws.Range("H9").Value = 13.81830698560672
$ws.Range("I9").Value = 22.93717509711494
$ws.Range("L9").Value = 9.341840360411052
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 21.71054551788774

# Row 10
$ws.Range("C10").Value = 10.62699022019643
$ws.Range("D10").Value = 3.84629491770482
$ws.Range("E10").Value = 12.97141709788208
$ws.Range("F10").Value = 25.77235484085012
$ws.Range("G10").Value = 33.05458916311118
$ws.Range("H10").Value = 13.91325508847762
$ws.Range("I10").Value = 23.3008192852084
$ws.Range("L10").Value = 9.359385048239249
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 22.09180496302788

# Row 11
$ws.Range("C11").Value = 10.66228121176306
$ws.Range("D11").Value = 3.847802294688922
$ws.Range("E11").Value = 13.01971159451615
$ws.Range("F11").Value = 26.06697412088463
$ws.Range("G11").Value = 33.62988588454041
$ws.Range("H11").Value = 13.95927001000017
$ws.Range("I11").Value = 23.47031043443026
$ws.Range("L11").Value = 9.369114970458936
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 22.26874091314627

# Row 12
$ws.Range("C12").Value = 10.67594794756019
$ws.Range("D12").Value = 3.848360398001512
$ws.Range("E12").Value = 13.0383826610261
$ws.Range("F12").Value = 26.17866826080337
$ws.Range("G12").Value = 33.84646349813296
$ws.Range("H12").Value = 13.97709061067313
$ws.Range("I12").Value = 23.53502473213233
$ws.Range("L12").Value = 9.373049282135611
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 22.33618803746311

# Row 13
$ws.Range("C13").Value = 10.67299123339524
$ws.Range("D13").Value = 3.848240766474188
$ws.Range("E13").Value = 13.03434466247171
$ws.Range("F13").Value = 26.154608986113
$ws.Range("G13").Value = 33.79987991998842
$ws.Range("H13").Value = 13.97323519235046
$ws.Range("I13").Value = 23.52106455949439
$ws.Range("L13").Value = 9.372190876534379
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 22.32164321214134

# Row 14
$ws.Range("C14").Value = 10.663399560055
$ws.Range("D14").Value = 3.847848465165371
$ws.Range("E14").Value = 13.02124006446528
$ws.Range("F14").Value = 26.0761613290333
$ws.Range("G14").Value = 33.64773069960471
$ws.Range("H14").Value = 13.96072825039572
$ws.Range("I14").Value = 23.47562417072503
$ws.Range("L14").Value = 9.369433652557962
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 22.2742812366249

# Row 15
$ws.Range("C15").Value = 10.65756357761039
$ws.Range("D15").Value = 3.847606512817502
$ws.Range("E15").Value = 13.01326265647562
$ws.Range("F15").Value = 26.02812328403841
$ws.Range("G15").Value = 33.55436215727666
$ws.Range("H15").Value = 13.95311861412575
$ws.Range("I15").Value = 23.44785826721328
$ws.Range("L15").Value = 9.367777251676774
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 22.24532692856539

# Row 16
$ws.Range("C16").Value = 10.62472667925071
$ws.Range("D16").Value = 3.846194631687085
$ws.Range("E16").Value = 12.96831515068674
$ws.Range("F16").Value = 25.75312477416261
$ws.Range("G16").Value = 33.01682629264795
$ws.Range("H16").Value = 13.91030385047433
$ws.Range("I16").Value = 23.28981995863531
$ws.Range("L16").Value = 9.358784226635981
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 22.08030725983653

# Row 17
$ws.Range("C17").Value = 10.60513012757246
$ws.Range("D17").Value = 3.845305910785564
$ws.Range("E17").Value = 12.9414355183911
$ws.Range("F17").Value = 25.58476974013747
$ws.Range("G17").Value = 32.68504285504687
$ws.Range("H17").Value = 13.88475429673296
$ws.Range("I17").Value = 23.19387278324213
$ws.Range("L17").Value = 9.35371407584595
$ws.Range("N17").Value = 20.9381421901699
$ws.Range("O17").Value = 21.97992825560852

# Row 18
$ws.Range("C18").Value = 10.59406251079566
$ws.Range("D18").Value = 3.844786461840909
$ws.Range("E18").Value = 12.92623370685793
$ws.Range("F18").Value = 25.48810033995069
$ws.Range("G18").Value = 32.49354347185814
$ws.Range("H18").Value = 13.87032518876805
$ws.Range("I18").Value = 23.1390726165462
$ws.Range("L18").Value = 9.350962573605424
$ws.Range("N18").Value = 20.79000725568364
$ws.Range("O18").Value = 21.92252576625879

# Row 19
$ws.Range("C19").Value = 10.59035043845728
$ws.Range("D19").Value = 3.844609164858994
$ws.Range("E19").Value = 12.92113138313668
$ws.Range("F19").Value = 25.45540136201334
$ws.Range("G19").Value = 32.42859811871119
$ws.Range("H19").Value = 13.86548577975245
$ws.Range("I19").Value = 23.1205861446342
$ws.Range("L19").Value = 9.350059299437921
$ws.Range("N19").Value = 20.73962067985786
$ws.Range("O19").Value = 21.90314919282143

# Row 20
$ws.Range("C20").Value = 10.60719517618117
$ws.Range("D20").Value = 3.845401373986273
$ws.Range("E20").Value = 12.94427021076845
$ws.Range("F20").Value = 25.60267527220227
$ws.Range("G20").Value = 32.72043249904289
$ws.Range("H20").Value = 13.88744659862088
$ws.Range("I20").Value = 23.20404694578725
$ws.Range("L20").Value = 9.354236766142472
$ws.Range("N20").Value = 20.96544799484618
$ws.Range("O20").Value = 21.99057973980244

# Row 21
$ws.Range("C21").Value = 10.66620871351054
$ws.Range("D21").Value = 3.847964038868174
$ws.Range("E21").Value = 13.02507889982575
$ws.Range("F21").Value = 26.09920070939486
$ws.Range("G21").Value = 33.69245701387599
$ws.Range("H21").Value = 13.96439118574074
$ws.Range("I21").Value = 23.48895711292354
$ws.Range("L21").Value = 9.370236750874012
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("O21").Value = 22.28818097415044

# Row 22
$ws.Range("C22").Value = 10.7065384941542
$ws.Range("D22").Value = 3.849564794649635
$ws.Range("E22").Value = 13.08011881599129
$ws.Range("F22").Value = 26.42440601491071
$ws.Range("G22").Value = 34.32020302926971
$ws.Range("H22").Value = 14.01698051735743
$ws.Range("I22").Value = 23.67823623155879
$ws.Range("L22").Value = 9.3821487773586
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 22.48524919820985

# Row 23
$ws.Range("C23").Value = 10.68485529430153
$ws.Range("D23").Value = 3.848717238062495
$ws.Range("E23").Value = 13.05054305136357
$ws.Range("F23").Value = 26.25081040271152
$ws.Range("G23").Value = 33.9859258456891
$ws.Range("H23").Value = 13.98870551089878
$ws.Range("I23").Value = 23.57695091822147
$ws.Range("L23").Value = 9.375658567017156
$ws.Range("N23").Value = 21.92877110912574
$ws.Range("O23").Value = 22.37985435241043

# Row 24
$ws.Range("C24").Value = 10.60626094819111
$ws.Range("D24").Value = 3.845358241561974
$ws.Range("E24").Value = 12.94298786167005
$ws.Range("F24").Value = 25.59457979821085
$ws.Range("G24").Value = 32.70443516708613
$ws.Range("H24").Value = 13.88622859925318
$ws.Range("I24").Value = 23.1994460778512
$ws.Range("L24").Value = 9.353999948776782
$ws.Range("N24").Value = 20.95310750188672
$ws.Range("O24").Value = 21.98576324384719

# Row 25
$ws.Range("C25").Value = 10.52930253894799
$ws.Range("D25").Value = 3.841418667015914
$ws.Range("E25").Value = 12.83689861319018
$ws.Range("F25").Value = 24.89346810154559
$ws.Range("G25").Value = 31.29776813676083
$ws.Range("H25").Value = 13.78590716680946
$ws.Range("I25").Value = 22.80721247484083
$ws.Range("L25").Value = 9.341840360411052
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 21.57362437670867

